$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.NumberFormat = "General"
    $c.Style = "Normal"
}

Set-TextValue "D2" "50.745.47"
Set-TextValue "D3" "2.924.31"
$ws.Range("E3").Value = "  -0.86%  "
Set-TextValue "D4" "1.00"
$ws.Range("E4").Value = "  +0.03%  "
Set-TextValue "D5" "375.72"
$ws.Range("E5").Value = "  -0.80%  "
Set-TextValue "D6" "99.50"
$ws.Range("E6").Value = "  -2.46%  "
$ws.Range("E7").Value = "  -0.32%  "
$ws.Range("E8").Value = "  -0.01%  "
Set-TextValue "D10" "35.63"
$ws.Range("E10").Value = "  -1.82%  "
$ws.Range("E11").Value = "  -0.47%  "
Set-TextValue "D12" "0.0850"
$ws.Range("E12").Value = "  +1.59%  "
Set-TextValue "D13" "3.385.71"
$ws.Range("E13").Value = "  -0.71%  "
Set-TextValue "D14" "18.02"
$ws.Range("E14").Value = "  +0.74%  "
Set-TextValue "D15" "7.59"
$ws.Range("E15").Value = "  +2.56%  "
Set-TextValue "D16" "12.00"
$ws.Range("E16").Value = "  +65.08%  "
Set-TextValue "D17" "2.930.48"
$ws.Range("E17").Value = "  -0.37%  "
Set-TextValue "D18" "0.988"
$ws.Range("E18").Value = "  +0.91%  "
Set-TextValue "D19" "50.722.08"
$ws.Range("E20").Value = "  -6.33%  "
$ws.Range("E21").Value = "  -1.31%  "
$ws.Range("E22").Value = "  -0.50%  "
Set-TextValue "D23" "69.36"
$ws.Range("E23").Value = "  +1.42%  "
Set-TextValue "D24" "265.51"
$ws.Range("E24").Value = "  +1.87%  "
Set-TextValue "D25" "3.16"
$ws.Range("E25").Value = "  +11.33%  "
$ws.Range("E26").Value = "  -3.88%  "
$ws.Range("E27").Value = "  -0.02%  "
Set-TextValue "D28" "7.04"
$ws.Range("E28").Value = "  -7.31%  "
$ws.Range("E29").Value = "  -1.11%  "
$ws.Range("E30").Value = "  -2.03%  "
$ws.Range("E31").Value = "  -4.19%  "
$ws.Range("E33").Value = "  -0.16%  "
$ws.Range("E34").Value = "  +0.12%  "
Set-TextValue "D35" "33.07"
$ws.Range("E35").Value = "  -1.12%  "
$ws.Range("E36").Value = "  -2.74%  "
$ws.Range("E37").Value = "  -0.02%  "
Set-TextValue "D38" "3.06"
$ws.Range("E38").Value = "  +3.61%  "
$ws.Range("E39").Value = "  +0.57%  "
Set-TextValue "D40" "16.29"
$ws.Range("E40").Value = "  -3.34%  "
$ws.Range("E41").Value = "  +1.15%  "
Set-TextValue "D42" "123.43"
$ws.Range("E42").Value = "  +1.61%  "
$ws.Range("E43").Value = "  -4.24%  "
Set-TextValue "D44" "20.88"
$ws.Range("E44").Value = "  -0.44%  "
$ws.Range("E45").Value = "  +6.12%  "
Set-TextValue "D46" "2.02"
$ws.Range("E46").Value = "  -1.48%  "
Set-TextValue "D48" "1.997.05"
$ws.Range("E48").Value = "  -0.05%  "
$ws.Range("E49").Value = "  -5.56%  "
$ws.Range("E50").Value = "  -5.15%  "
$ws.Range("E51").Value = "  +4.03%  "
